# Data cleanup: fix inconsistent player_id values for LeBron James rows.
# Rows 4-23 (column C, player_id) incorrectly contain 3463; they should be
# 3462, matching the rest of LeBron James's rows (rows 2-3, 24 onward).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4:C23").Value = 3462
